$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '40.115.60'
$ws.Range("E2").Value = '  +2.20%  '

$ws.Range("D3").Value = '2.239.00'
$ws.Range("E3").Value = '  +1.54%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '293.20'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.77%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '87.57'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.81%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.517'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.58%  '

$ws.Range("E8").Value = '  -0.17%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.474'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.90%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '31.41'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.53%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0793'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.65%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.14'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.71%  '

$ws.Range("E13").Value = '  +1.42%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.43'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.89%  '

$ws.Range("D15").Value = '2.577.92'
$ws.Range("E15").Value = '  +1.47%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.13'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.34%  '

$ws.Range("D17").Value = '2.224.32'
$ws.Range("E17").Value = '  +1.31%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.734'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.19%  '

$ws.Range("D19").Value = '39.977.58'
$ws.Range("E19").Value = '  +2.18%  '

$ws.Range("D20").Value = '0.0₃0889'
$ws.Range("E20").Value = '  +2.22%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.24'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +9.29%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.84'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.86%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.64'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.64%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '236.54'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.82%  '

$ws.Range("E25").Value = '  +0.00%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.48'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.40%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.86'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.76%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.96'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.04%  '

$ws.Range("E29").Value = '  +3.14%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.35'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.81%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '33.52'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.88%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '152.40'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.04%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.997'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.26%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.95'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.56%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0724'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.93%  '

$ws.Range("E36").Value = '  +2.66%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '16.46'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +8.72%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.84'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.27%  '

$ws.Range("E39").Value = '  +2.94%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.100'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.77%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.74'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.55%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.85'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.56%  '

$ws.Range("D43").Value = '2.077.26'
$ws.Range("E43").Value = '  +9.17%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '18.29'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +15.52%  '

$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0270'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.93%  '

$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.11'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.97%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.94'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +12.00%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.61'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.08%  '

$ws.Range("D49").Value = '2.437.72'
$ws.Range("E49").Value = '  +1.19%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '72.01'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.45%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '89.58'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.62%  '
